$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Round")
$ws.Activate()

$ws.Range("A4").Value = "2020100,2020101,2020102"
$ws.Range("A5").Value = "2021103,2021104,2022105"
$ws.Range("B5").Value = "3,3,1"

$ws.Range("I14").Select()
